$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("Y25").Value = 0
$ws.Range("AB28").Value = 0
$ws.Range("AE31").Value = 0
$ws.Range("AF32").Value = 0
$ws.Range("AG33").Value = 0
